$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Project Overview")

# The "Product Backlog" list (columns I/J, zig-zagging row to row) gets two
# new entries inserted right before "Meal Planner/Grocery list":
#   - "Modify Recipe"
#   - "Step re-order (create easier process)"
# Every entry from "Meal Planner/Grocery list" onward therefore shifts down
# by one slot in the list (I25->I27, I26->I28, J27->J29, I28->I30, J29->J31,
# I30->I32).

# New slots: the two newly-added backlog items.
$ws.Range("I25").Value = "Modify Recipe"
$ws.Range("J26").Value = "Step re-order (create easier process)"

# Old I26 slot is vacated (its content moved to I28 below).
$ws.Range("I26").Value = $null

# Previously existing items, shifted down by one slot.
$ws.Range("I27").Value = "Meal Planner/Grocery list"
$ws.Range("I28").Value = "Account settings"

# Old J27 slot is vacated (its content moved to J29 below).
$ws.Range("J27").Value = $null

$ws.Range("J29").Value = "Themes/layouts"
$ws.Range("I30").Value = "Master Ingredient List"
$ws.Range("J31").Value = "Allow user to add to"
$ws.Range("I32").Value = "IOS"

# Match the saved view state: scrolled down a bit, with J27 selected.
$ws.Activate()
$ws.Range("A7").Select()
$excel.ActiveWindow.ScrollRow = 7
$ws.Range("J27").Select()
